$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 22: a "Sequence generator" reset entry, duplicating
# --- the Field/Type/Generate of row 21 but with a new Pattern value
# --- that implements the SEQ(...) reset.

$ws.Range("A22").Value = "Sequence generator"
$ws.Range("B22").Value = "Sequence"
$ws.Range("C22").Value = "y"
$ws.Range("D22").Value = "[SEQ(4,1,RESTART)]"

# Row height to match the rest of the table
$ws.Rows("22").RowHeight = 15.75

# Borders matching the look of row 21 (thin black/automatic edges)
# A22: left + right thin black
$ws.Range("A22").Borders.Item(7).LineStyle = 1
$ws.Range("A22").Borders.Item(7).Color = 0
$ws.Range("A22").Borders.Item(10).LineStyle = 1
$ws.Range("A22").Borders.Item(10).Color = 0

# B22: right thin black
$ws.Range("B22").Borders.Item(10).LineStyle = 1
$ws.Range("B22").Borders.Item(10).Color = 0

# C22: left thin black, right thin automatic
$ws.Range("C22").Borders.Item(7).LineStyle = 1
$ws.Range("C22").Borders.Item(7).Color = 0
$ws.Range("C22").Borders.Item(10).LineStyle = 1

# D22: right thin automatic, text number format (matches pattern column)
$ws.Range("D22").Borders.Item(10).LineStyle = 1
$ws.Range("D22").NumberFormat = "@"

# --- Update the view state: scroll position and active selection
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G16").Select()
